$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new value pairs (price/volume refresh from data source)
$updates = @(
    @{ Cell = 'D2'; Value = '44.186.62' }
    @{ Cell = 'E2'; Value = '  +5.20%  ' }
    @{ Cell = 'D3'; Value = '2.265.55' }
    @{ Cell = 'E3'; Value = '  +2.53%  ' }
    @{ Cell = 'E4'; Value = '  +0.12%  ' }
    @{ Cell = 'D5'; Value = '230.58' }
    @{ Cell = 'E5'; Value = '  +0.13%  ' }
    @{ Cell = 'D6'; Value = '0.633' }
    @{ Cell = 'E6'; Value = '  +2.60%  ' }
    @{ Cell = 'D7'; Value = '63.70' }
    @{ Cell = 'E7'; Value = '  +4.98%  ' }
    @{ Cell = 'E8'; Value = '  +0.01%  ' }
    @{ Cell = 'D9'; Value = '0.447' }
    @{ Cell = 'E9'; Value = '  +11.28%  ' }
    @{ Cell = 'D10'; Value = '0.103' }
    @{ Cell = 'E10'; Value = '  +14.92%  ' }
    @{ Cell = 'D11'; Value = '56.84' }
    @{ Cell = 'E11'; Value = '  -0.67%  ' }
    @{ Cell = 'D12'; Value = '26.35' }
    @{ Cell = 'E12'; Value = '  +19.55%  ' }
    @{ Cell = 'E13'; Value = '  +2.61%  ' }
    @{ Cell = 'D14'; Value = '2.598.10' }
    @{ Cell = 'E14'; Value = '  +2.35%  ' }
    @{ Cell = 'D15'; Value = '15.72' }
    @{ Cell = 'E15'; Value = '  +1.96%  ' }
    @{ Cell = 'D16'; Value = '6.08' }
    @{ Cell = 'E16'; Value = '  +9.07%  ' }
    @{ Cell = 'D18'; Value = '2.258.93' }
    @{ Cell = 'E18'; Value = '  +2.43%  ' }
    @{ Cell = 'D19'; Value = '43.982.35' }
    @{ Cell = 'E19'; Value = '  +5.14%  ' }
    @{ Cell = 'E20'; Value = '  +7.33%  ' }
    @{ Cell = 'D21'; Value = '73.72' }
    @{ Cell = 'E21'; Value = '  +2.09%  ' }
    @{ Cell = 'E22'; Value = '  -0.43%  ' }
    @{ Cell = 'D23'; Value = '255.44' }
    @{ Cell = 'E23'; Value = '  +5.32%  ' }
    @{ Cell = 'E24'; Value = '  +0.11%  ' }
    @{ Cell = 'D25'; Value = '2.43' }
    @{ Cell = 'E25'; Value = '  +3.47%  ' }
    @{ Cell = 'E26'; Value = '  -2.13%  ' }
    @{ Cell = 'B27'; Value = 'Cosmos' }
    @{ Cell = 'C27'; Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom' }
    @{ Cell = 'D27'; Value = '10.14' }
    @{ Cell = 'E27'; Value = '  +5.45%  ' }
    @{ Cell = 'B28'; Value = 'WEMIXToken' }
    @{ Cell = 'C28'; Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix' }
    @{ Cell = 'D28'; Value = '3.34' }
    @{ Cell = 'E28'; Value = '  +25.96%  ' }
    @{ Cell = 'D29'; Value = '172.04' }
    @{ Cell = 'E29'; Value = '  +1.87%  ' }
    @{ Cell = 'D30'; Value = '20.80' }
    @{ Cell = 'E30'; Value = '  +2.19%  ' }
    @{ Cell = 'E31'; Value = '  -1.57%  ' }
    @{ Cell = 'E32'; Value = '  -2.68%  ' }
    @{ Cell = 'E33'; Value = '  +3.11%  ' }
    @{ Cell = 'D34'; Value = '0.0679' }
    @{ Cell = 'E34'; Value = '  +4.98%  ' }
    @{ Cell = 'E35'; Value = '  +3.84%  ' }
    @{ Cell = 'E36'; Value = '  -1.42%  ' }
    @{ Cell = 'D37'; Value = '3.84' }
    @{ Cell = 'E37'; Value = '  +8.61%  ' }
    @{ Cell = 'D38'; Value = '6.72' }
    @{ Cell = 'E38'; Value = '  +7.11%  ' }
    @{ Cell = 'E39'; Value = '  -0.35%  ' }
    @{ Cell = 'E40'; Value = '  +5.23%  ' }
    @{ Cell = 'E41'; Value = '  -0.01%  ' }
    @{ Cell = 'B42'; Value = 'InjectiveProtocol' }
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj' }
    @{ Cell = 'D42'; Value = '17.53' }
    @{ Cell = 'E42'; Value = '  +9.22%  ' }
    @{ Cell = 'B43'; Value = 'FraxShare' }
    @{ Cell = 'C43'; Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs' }
    @{ Cell = 'D43'; Value = '8.32' }
    @{ Cell = 'E43'; Value = '  -2.96%  ' }
    @{ Cell = 'D44'; Value = '0.0966' }
    @{ Cell = 'E44'; Value = '  +1.44%  ' }
    @{ Cell = 'D45'; Value = '97.84' }
    @{ Cell = 'E45'; Value = '  +1.27%  ' }
    @{ Cell = 'B46'; Value = 'FTXToken' }
    @{ Cell = 'C46'; Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt' }
    @{ Cell = 'D46'; Value = '4.39' }
    @{ Cell = 'E46'; Value = '  +0.40%  ' }
    @{ Cell = 'B47'; Value = 'TrustWalletToken' }
    @{ Cell = 'C47'; Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt' }
    @{ Cell = 'D47'; Value = '1.19' }
    @{ Cell = 'E47'; Value = '  +0.18%  ' }
    @{ Cell = 'B48'; Value = 'TerraClassic' }
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc' }
    @{ Cell = 'D48'; Value = '0.000211' }
    @{ Cell = 'E48'; Value = '  -6.36%  ' }
    @{ Cell = 'D49'; Value = '10.06' }
    @{ Cell = 'E49'; Value = '  +18.11%  ' }
    @{ Cell = 'D50'; Value = '1.447.11' }
    @{ Cell = 'E50'; Value = '  -0.63%  ' }
    @{ Cell = 'D51'; Value = '2.30' }
    @{ Cell = 'E51'; Value = '  +4.16%  ' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    # Force text storage so numeric-looking strings (e.g. "230.58") are not
    # coerced into real numbers, matching the source data which is text.
    $rng.NumberFormat = '@'
    $rng.Value = $u.Value
    $rng.Style = 'Normal'
}
